$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws.Range("D2").Value = 2879.62
$ws.Range("E2").Value = -2879.62

$ws.Range("D3").Value = 943.51
$ws.Range("E3").Value = 12779.83
$ws.Range("F3").Value = 0.06875221338245646

$ws.Range("D4").Value = 3823.13
$ws.Range("E4").Value = 9900.209999999999
$ws.Range("F4").Value = 0.2785859710536939
